# BCHEP-542: contractor import process refactor
# - Replace plain-text contractor email addresses in column A (Username)
#   with anonymized "test.exampleN@fakeemail.com" display values that are
#   turned into mailto: hyperlinks (pointing at the real address), and
# - Reset the sheet view (drop the stale scrolled/selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the plain-text contractor e-mail addresses that used to live in
# column A (Username) with anonymized "test.exampleN@fakeemail.com" display
# values. Processed in this specific order so the newly-interned shared
# strings land in the same order as the source workbook (row3, row4, row2,
# row5, row6).
$ws.Range("A3").Value = "test.example2@fakeemail.com"
$ws.Range("A4").Value = "test.example3@fakeemail.com"
$ws.Range("A2").Value = "test.example1@fakeemail.com"
$ws.Range("A5").Value = "test.example4@fakeemail.com"
$ws.Range("A6").Value = "test.example5@fakeemail.com"

# Turn each of those cells into a mailto: hyperlink pointing at the real
# contractor address, added in sheet order (A2..A6).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:Pridegasheating@yahoo.ca")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:marklandelectric@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:dan@authenticinstallations.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:valleysidehvac@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:Devon@harbourhazmat.ca")

# Restore the view: select A7 and scroll back so AK1 is no longer pinned at
# the top-left.
$ws.Range("A7").Select()
